$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (256-269): column A is the date serial, B/C/D are numeric values.
$data = @(
    @(256, 44330, 0, 0, 0),
    @(257, 44331, 0, 0, 0),
    @(258, 44332, 0, 0, 0),
    @(259, 44333, 0, 0, 0),
    @(260, 44334, 0, 0, 0),
    @(261, 44335, 0, 0, 0),
    @(262, 44336, 2, 2, 93.41429238673517),
    @(263, 44337, 0, 2, 93.41429238673517),
    @(264, 44338, 0, 2, 93.41429238673517),
    @(265, 44339, 0, 2, 93.41429238673517),
    @(266, 44340, 0, 2, 93.41429238673517),
    @(267, 44341, 0, 2, 93.41429238673517),
    @(268, 44342, 0, 2, 93.41429238673517),
    @(269, 44343, 0, 0, 0)
)

# Reference cell (A255) whose style (bold, border, date number format) must be
# replicated onto the new column-A date cells.
$styleSource = $ws.Range("A255")

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value2 = $r[1]
    $ws.Cells.Item($row, 2).Value2 = $r[2]
    $ws.Cells.Item($row, 3).Value2 = $r[3]
    $ws.Cells.Item($row, 4).Value2 = $r[4]

    $styleSource.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
